$wb = $excel.ActiveWorkbook

# Rows (1-based sheet rows) holding the 8 redundant metrics: MBRAE, UMBRAE,
# STDAPE, RMSPE, MRE, MRAE, MDRAE, GMRAE
$rowsToDelete = 12, 13, 23, 24, 29, 31, 32, 33

foreach ($ws in $wb.Worksheets) {
    foreach ($r in ($rowsToDelete | Sort-Object -Descending)) {
        $ws.Rows.Item($r).EntireRow.Delete()
    }

    # Renumber the ID column (A) sequentially for the remaining data rows
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}
